$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cck"
$ws.Range("C2").Value = "Cckar"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.862015666666667
$ws.Range("H2").Value = 5.586047
$ws.Range("I2").Value = 0.399282427039184
$ws.Range("J2").Value = 0.3992824270391842
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1497266666666667
$ws.Range("N2").Value = 0.44918
$ws.Range("O2").Value = 0.4607938481487409
$ws.Range("P2").Value = 0.4607938481487409
$ws.Range("Q2").Value = 0.2787933990511111
$ws.Range("R2").Value = 2.50914059146
$ws.Range("S2").Value = 0.1839868860535545
$ws.Range("T2").Value = 0.1839868860535545

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cck"
$ws.Range("C3").Value = "Cckar"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.862015666666667
$ws.Range("H3").Value = 5.586047
$ws.Range("I3").Value = 0.399282427039184
$ws.Range("J3").Value = 0.3992824270391842
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1752053333333333
$ws.Range("N3").Value = 0.525616
$ws.Range("O3").Value = 0.5392061518512591
$ws.Range("P3").Value = 0.5392061518512591
$ws.Range("Q3").Value = 0.3262350755502222
$ws.Range("R3").Value = 2.936115679952
$ws.Range("S3").Value = 0.2152955409856296
$ws.Range("T3").Value = 0.2152955409856296

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cck"
$ws.Range("C4").Value = "Cckar"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.017090333333333
$ws.Range("H4").Value = 6.051271
$ws.Range("I4").Value = 0.4325359545939786
$ws.Range("J4").Value = 0.4325359545939787
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1497266666666667
$ws.Range("N4").Value = 0.44918
$ws.Range("O4").Value = 0.4607938481487409
$ws.Range("P4").Value = 0.4607938481487409
$ws.Range("Q4").Value = 0.3020122119755556
$ws.Range("R4").Value = 2.71810990778
$ws.Range("S4").Value = 0.1993099069800485
$ws.Range("T4").Value = 0.1993099069800485

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cck"
$ws.Range("C5").Value = "Cckar"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.017090333333333
$ws.Range("H5").Value = 6.051271
$ws.Range("I5").Value = 0.4325359545939786
$ws.Range("J5").Value = 0.4325359545939787
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1752053333333333
$ws.Range("N5").Value = 0.525616
$ws.Range("O5").Value = 0.5392061518512591
$ws.Range("P5").Value = 0.5392061518512591
$ws.Range("Q5").Value = 0.3534049842151111
$ws.Range("R5").Value = 3.180644857936
$ws.Range("S5").Value = 0.2332260476139301
$ws.Range("T5").Value = 0.2332260476139302

$ws.Range("A6").Value = "Inflammatory-Mac"
$ws.Range("B6").Value = "Cck"
$ws.Range("C6").Value = "Cckar"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1819486666666667
$ws.Range("H6").Value = 0.5458460000000001
$ws.Range("I6").Value = 0.0390162695855639
$ws.Range("J6").Value = 0.03901626958556391
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1497266666666667
$ws.Range("N6").Value = 0.44918
$ws.Range("O6").Value = 0.4607938481487409
$ws.Range("P6").Value = 0.4607938481487409
$ws.Range("Q6").Value = 0.02724256736444445
$ws.Range("R6").Value = 0.24518310628
$ws.Range("S6").Value = 0.01797845700274067
$ws.Range("T6").Value = 0.01797845700274068

$ws.Range("A7").Value = "Inflammatory-Mac"
$ws.Range("B7").Value = "Cck"
$ws.Range("C7").Value = "Cckar"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1819486666666667
$ws.Range("H7").Value = 0.5458460000000001
$ws.Range("I7").Value = 0.0390162695855639
$ws.Range("J7").Value = 0.03901626958556391
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.1752053333333333
$ws.Range("N7").Value = 0.525616
$ws.Range("O7").Value = 0.5392061518512591
$ws.Range("P7").Value = 0.5392061518512591
$ws.Range("Q7").Value = 0.03187837679288889
$ws.Range("R7").Value = 0.286905391136
$ws.Range("S7").Value = 0.02103781258282323
$ws.Range("T7").Value = 0.02103781258282323

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Cck"
$ws.Range("C8").Value = "Cckar"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.6023503333333334
$ws.Range("H8").Value = 1.807051
$ws.Range("I8").Value = 0.1291653487812732
$ws.Range("J8").Value = 0.1291653487812732
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.1497266666666667
$ws.Range("N8").Value = 0.44918
$ws.Range("O8").Value = 0.4607938481487409
$ws.Range("P8").Value = 0.4607938481487409
$ws.Range("Q8").Value = 0.09018790757555557
$ws.Range("R8").Value = 0.8116911681800001
$ws.Range("S8").Value = 0.05951859811239716
$ws.Range("T8").Value = 0.05951859811239716

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Cck"
$ws.Range("C9").Value = "Cckar"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.6023503333333334
$ws.Range("H9").Value = 1.807051
$ws.Range("I9").Value = 0.1291653487812732
$ws.Range("J9").Value = 0.1291653487812732
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1752053333333333
$ws.Range("N9").Value = 0.525616
$ws.Range("O9").Value = 0.5392061518512591
$ws.Range("P9").Value = 0.5392061518512591
$ws.Range("Q9").Value = 0.1055349909351111
$ws.Range("R9").Value = 0.9498149184160001
$ws.Range("S9").Value = 0.06964675066887605
$ws.Range("T9").Value = 0.06964675066887605

$ws.Rows("10:10").Delete()